$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 263, shifting existing rows 263:354 down to 264:355.
$ws.Rows.Item(263).Insert()

# Populate the newly inserted row 263 with the new data record.
$ws.Range("A263").Value = 4
$ws.Range("B263").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C263").Value = "Los Lagos"
$ws.Range("D263").Value = 44809
$ws.Range("E263").Value = 10
$ws.Range("F263").Value = 100112045
$ws.Range("G263").Value = "Zapallo"
$ws.Range("H263").Value = "Paine"
$ws.Range("I263").Value = "1a (guarda)"
$ws.Range("J263").Value = 500
$ws.Range("K263").Value = 500
$ws.Range("L263").Value = 650
$ws.Range("M263").Value = 575
$ws.Range("N263").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O263").Value = "Región de O'Higgins"
$ws.Range("P263").Value = 575
$ws.Range("Q263").Value = 1
$ws.Range("R263").Value = "Hortaliza"
